$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the File Name column (C) values to the new naming scheme.
$ws.Range("C2").Value = "import_character.xml"
$ws.Range("C3").Value = "import_audio.xml"
$ws.Range("C4").Value = "import_character.xml"
$ws.Range("C5").Value = "import_book_chapter.xml"
$ws.Range("C6").Value = "import_book.xml"
$ws.Range("C7").Value = "import_image.xml"
$ws.Range("C8").Value = "import_event.xml"
$ws.Range("C9").Value = "import_material.xml"
$ws.Range("C10").Value = "import_location.xml"
$ws.Range("C11").Value = "import_archive.xml"

# Add a new row 12 for O_011
$ws.Range("A12").Value = "O_011"
$ws.Range("C12").Value = "import_documentation.xml"
$ws.Range("D12").Value = "data/XML/"

# Update Description column (B) for all data rows to "to be described"
$ws.Range("B2:B12").Value = "to be described"

# Update the selection to mirror the saved cursor position
$ws.Range("B17").Select()
